$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: new SIQ entry
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "user can access the elevator by his username and password "
$ws.Range("E4").Value = "we can use keypad and each number in keypad has array of few character  as nokia mobile 6630"
$ws.Range("D4").Value = "how can user enter his username?"
$ws.Range("C4").Value = "elevator has lock system so must enter userID to access elevator"
$ws.Range("F4").Value = "22/1/2020"
$ws.Range("G4").Value = "24/1/2020"

# Row 3: Comments column gets "not answered"
$ws.Range("J3").Value = "not answered"
$ws.Range("J4").Value = "not answered"

# Update the active selection to I10
[void]$ws.Range("I10").Select()
